# Re-point the displayed "URL" value on the Global sheet from the old
# custodian/test host (.111) to the new one (.107). The hyperlink
# relationship itself (rId1 -> http://192.168.168.111/) is left
# untouched - only the cell's visible text (the shared string) changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "http://192.168.168.107/"

# Scroll the view and move the selection from Y2 (scrolled to column R)
# to D3 (scrolled back to column D).
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1

$ws.Range("D3").Select()
